$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.725.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.82%  "

$ws.Range("D3").Value = "'2.233.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'231.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").Value = "'0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("D7").Value = "'60.95"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "'0.0914"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "'2.569.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").Value = "'15.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").Value = "'22.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "

$ws.Range("E15").Value = "  +2.13%  "

$ws.Range("D16").Value = "'0.801"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").Value = "'2.237.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").Value = "'42.550.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.75%  "

$ws.Range("D19").Value = "'0.0₃0939"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.35%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'72.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("D22").Value = "'245.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").Value = "'2.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.68%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  +6.67%  "

$ws.Range("D26").Value = "'9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'169.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.143"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("E30").Value = "  +3.40%  "

$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("D33").Value = "'5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("D34").Value = "'4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "

$ws.Range("D35").Value = "'0.0654"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("D36").Value = "'6.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "

$ws.Range("D38").Value = "'3.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "

$ws.Range("E39").Value = "  +3.93%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").Value = "'0.000224"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.59%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0967"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.66%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'97.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.08%  "

$ws.Range("D47").Value = "'1.459.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").Value = "'16.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").Value = "'2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.93%  "
